# kaizen: añadiendo oficiales y no oficiales
# Adds two new header columns ("Oficiales" / "No. Oficiales") to both
# worksheets (AX ENGOMADO -> P1:Q1, AX URDIDO -> N1:O1), and flips which
# sheet/selection is active.

$wb = $excel.ActiveWorkbook
$wsEngomado = $wb.Worksheets.Item(1)   # "AX ENGOMADO"
$wsUrdido   = $wb.Worksheets.Item(2)   # "AX URDIDO"

# --- AX ENGOMADO (sheet1): add P1 "Oficiales", Q1 "No. Oficiales" ---
$wsEngomado.Range("P1").Value = "Oficiales"
$wsEngomado.Range("Q1").Value = "No. Oficiales"

$rngEng = $wsEngomado.Range("P1:Q1")
$rngEng.HorizontalAlignment = -4108
$rngEng.VerticalAlignment = -4108
$rngEng.WrapText = $true

# --- AX URDIDO (sheet2): add N1 "Oficiales", O1 "No. Oficiales" ---
$wsUrdido.Range("N1").Value = "Oficiales"
$wsUrdido.Range("O1").Value = "No. Oficiales"

$rngUrd = $wsUrdido.Range("N1:O1")
$rngUrd.HorizontalAlignment = -4108
$rngUrd.VerticalAlignment = -4108
$rngUrd.WrapText = $true

# --- Selection / active sheet: AX URDIDO first (so its own selection is
#     recorded), then AX ENGOMADO last, becoming the active tab. ---
$wsUrdido.Range("N1:O1").Select()
$wsEngomado.Range("P1:Q1").Select()
